# Updates cryptos list values (price/volume) per diff; also swaps rows 48/49 (Bittensor <-> USDe)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = "'61.574.06"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Formula = "'2.881.83"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Formula = "'  -2.36%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Formula = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Formula = "'  -0.04%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Formula = "'566.12"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Formula = "'  -4.54%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Formula = "'142.36"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Formula = "'  -3.64%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Formula = "'  +0.05%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Formula = "'  -1.10%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Formula = "'2.876.79"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Formula = "'  -2.44%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Formula = "'6.85"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Formula = "'  -2.77%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Formula = "'  -2.55%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Formula = "'  -2.16%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Formula = "'0.0000230"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Formula = "'  -1.73%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Formula = "'31.57"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Formula = "'  -2.95%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D16').Formula = "'3.358.08"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Formula = "'  -2.34%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Formula = "'61.518.27"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Formula = "'  -2.07%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Formula = "'2.886.34"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Formula = "'  -2.07%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D20').Formula = "'429.62"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Formula = "'  -1.89%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Formula = "'12.98"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Formula = "'  -3.10%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Formula = "'  -2.18%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Formula = "'  -3.12%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Formula = "'  -2.36%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Formula = "'11.79"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Formula = "'  -0.06%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Formula = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Formula = "'  -0.05%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Formula = "'9.91"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Formula = "'  -12.02%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Formula = "'  -6.14%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Formula = "'  +6.40%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Formula = "'  -3.99%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Formula = "'  -4.88%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E33').Formula = "'  +0.00%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Formula = "'0.105"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Formula = "'  -2.37%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Formula = "'25.37"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Formula = "'  -3.75%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Formula = "'0.952"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Formula = "'  -3.96%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Formula = "'  -4.81%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Formula = "'48.79"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Formula = "'  -1.58%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Formula = "'  -6.97%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Formula = "'  -6.15%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Formula = "'8.13"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Formula = "'  -3.33%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Formula = "'  -4.27%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Formula = "'39.40"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Formula = "'  -0.02%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Formula = "'  -5.26%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Formula = "'2.676.79"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Formula = "'  -0.46%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Formula = "'132.33"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Formula = "'  -2.07%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Formula = "'  -1.34%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Formula = "'USDe"
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Formula = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Formula = "'1.00"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Formula = "'  -0.01%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Formula = "'Bittensor"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Formula = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Formula = "'341.89"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Formula = "'  -4.40%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Formula = "'  -2.07%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Formula = "'21.36"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Formula = "'  -5.73%  "
$ws.Range('E51').Style = 'Normal'
